# Update the embedded OLS regression "Summary" report text on each of the
# 28 worksheets: the Date: line changes from "Thu, 02 Jan 2020" to
# "Sun, 05 Jan 2020", and the Time: line changes from "20:49:03" to a new
# re-run time ("21:22:43" on the first 9 sheets, "21:22:44" on the rest),
# matching a later re-execution of the backward-elimination script.

$wb = $excel.ActiveWorkbook

$oldDate = "Thu, 02 Jan 2020"
$newDate = "Sun, 05 Jan 2020"
$oldTime = "20:49:03"

$sheetCount = $wb.Worksheets.Count
Write-Host "Sheet count:" $sheetCount

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    if ($i -le 9) {
        $newTime = "21:22:43"
    } else {
        $newTime = "21:22:44"
    }

    $cell = $ws.Range("B2")
    $text = $cell.Value2
    $origRowHeight = $ws.Rows.Item(2).RowHeight

    if ($text -ne $null -and $text.Contains($oldDate)) {
        $text = $text.Replace($oldDate, $newDate)
        $text = $text.Replace($oldTime, $newTime)
        $cell.Value = $text
        # Writing the cell value causes the engine to re-run row autofit,
        # which can push the (already maxed-out) row height past Excel's
        # 409.5pt ceiling. Restore the original row height so only the
        # text content changes, matching the source edit.
        $ws.Rows.Item(2).RowHeight = $origRowHeight
        Write-Host "Updated sheet" $i $ws.Name
    } else {
        Write-Host "Skipped sheet" $i $ws.Name "(no match)"
    }
}
